$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G1").Value = "MappedValue"
$ws.Range("G2").Value = "a"
$ws.Range("G3").Value = "B"
$ws.Range("G4").Value = "b"
$ws.Range("G5").Value = "D"

$ws.Range("G6").Select()
